# weights.xlsx edit:
#  - rename Sheet1 -> "Weights"
#  - add a new "Stats" sheet after it
#  - insert 4 new precious-metal ETF tickers (IAU, SLV, PALL, PPLT) into the
#    weights table, which pushes PFE/SYY/NRZ/OHI into four new trailing rows
#  - update portfolio weights (col B) to the new allocation (sums to 1.0)
#  - populate Stats with the expected return / std dev / utility figures

$wb = $excel.ActiveWorkbook

# --- Weights sheet -----------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Weights"

$ws1.Range("A1").Value = "Ticker"
$ws1.Range("B1").Value = "Weight"

$tickers = @("TLT","GLTR","IAU","SLV","PALL","PPLT","FDX","WMT","AMZN","CORE","BLK","LMT","ORCL","NTRS","TSM","SJM","MDLZ","REGI","V","MSFT","JNJ","TPH","VIRT","AXP","BX","CNC","LDOS","MDT","MRK","NKE","PFE","SYY","NRZ","OHI")
$weights = @(0.1622,0.15,0,0,0,0,0.005,0.05,0.1,0.005,0.005,0.005,0.005,0.005,0.05,0.005,0.005,0.05,0.05,0.05,0.005,0.005,0.0428,0.005,0.05,0.005,0.05,0.005,0.05,0.05,0.005,0.005,0.005,0.015)

for ($i = 0; $i -lt $tickers.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 1).Value = $tickers[$i]
    $ws1.Cells.Item($row, 2).Value = $weights[$i]
}

# --- Stats sheet ---------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Stats"

# The figures are recorded as text (not numbers), so force column B to a
# text number-format before assigning, otherwise Excel auto-coerces a
# numeric-looking string into a real number.
$ws2.Range("B1:B3").NumberFormat = "@"

$ws2.Range("A1").Value = "Expected Annualized Return"
$ws2.Range("B1").Value = "0.2058"
$ws2.Range("A2").Value = "Standard Deviation"
$ws2.Range("B2").Value = "0.1583"
$ws2.Range("A3").Value = "Portfolio Utility"
$ws2.Range("B3").Value = "0.1933"

# Keep "Weights" as the active/visible tab, matching the target workbook.
$ws1.Activate()
